$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: _old -> _FV2410 (columns A-J = 1-10), _new -> _FV2504 (columns L-U = 12-21)
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_old$', '_FV2410')
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_new$', '_FV2504')
}

# Add table over the full data range (defaults to name "Table1" as the first table)
$range = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
if ($tbl.Name -ne "Table1") { $tbl.Name = "Table1" }

# Freeze top row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
